# Actualización SmartScore desde Streamlit (Harrison Driver)
# Adds a new response row (row 20) to the results sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 20

$ws.Cells.Item($row, 1).Value = "Harrison Driver_20251202_130401"
$ws.Cells.Item($row, 2).Value = "'"
$ws.Cells.Item($row, 3).Value = "Harrison Driver"
$ws.Cells.Item($row, 4).Value = 18
$ws.Cells.Item($row, 5).Value = "Male"
$ws.Cells.Item($row, 6).Value = "2025-12-02 13:04:01"

$pesos = @"
{
  "portion": 1.0,
  "diet": 1.0,
  "salt": 0.2,
  "fat": 0.2,
  "natural": 0.8,
  "convenience": 1.0,
  "price": 0.8
}
"@
$ws.Cells.Item($row, 7).Value = $pesos

$ws.Cells.Item($row, 8).Value = "Nongshim Neoguri Spicy Seafood"
$ws.Cells.Item($row, 9).Value = "'0.646"
$ws.Cells.Item($row, 10).Value = "Sabor a marisco, umami, picante equilibrado, buena textura, algo salado"

$ws.Cells.Item($row, 11).Value = "Nissin Chow Mein Teriyaki Beef"
$ws.Cells.Item($row, 12).Value = "'0.543"
$ws.Cells.Item($row, 13).Value = "Fácil de preparar, porción generosa, salsa suave, necesita mejoras, alto en grasa"

$ws.Cells.Item($row, 14).Value = "Nongshim Shin Ramyun"
$ws.Cells.Item($row, 15).Value = "'0.535"
$ws.Cells.Item($row, 16).Value = "Sabor intenso, picante, umami, fideos gruesos, muy alto en sodio"

$ws.Cells.Item($row, 17).Value = "Amy’s Macaroni & Cheese (frozen)"
$ws.Cells.Item($row, 18).Value = "'0.662"
$ws.Cells.Item($row, 19).Value = "Queso real, textura casera, sin conservadores, alto en grasa, algo caro"

$ws.Cells.Item($row, 20).Value = "Kraft Macaroni & Cheese Dinner"
$ws.Cells.Item($row, 21).Value = "'0.507"
$ws.Cells.Item($row, 22).Value = "Sabor nostálgico, clásico americano, fácil, no muy nutritivo, barato"

$ws.Cells.Item($row, 23).Value = "Annie’s Shells & White Cheddar"
$ws.Cells.Item($row, 24).Value = "'0.456"
$ws.Cells.Item($row, 25).Value = "Queso blanco real, sin colorantes, sabor casero, menos salado, buena para niños"

$ws.Cells.Item($row, 26).Value = "Wild Planet Wild Tuna Pasta Salad"
$ws.Cells.Item($row, 27).Value = "'0.720"
$ws.Cells.Item($row, 28).Value = "Sabor fresco, buena proteína, saludable, porción algo pequeña"

$ws.Cells.Item($row, 29).Value = "StarKist Chicken Creations (Chicken Salad)"
$ws.Cells.Item($row, 30).Value = "'0.498"
$ws.Cells.Item($row, 31).Value = "Portátil, saludable, fácil, buena textura, sabor suave"

$ws.Cells.Item($row, 32).Value = "Kitchens of India Variety Pack"
$ws.Cells.Item($row, 33).Value = "'0.472"
$ws.Cells.Item($row, 34).Value = "Sabor auténtico, variedad, vegetariano, necesita arroz o pan, buena calidad"
